$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Ciboulette (Vega Modelo de Temuco).
# It belongs right above the current row 326, so insert a blank row there -
# this pushes the existing rows 326-431 down to 327-432 (dimension grows to R432).
$ws.Rows(326).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(326, 1).Value = 10
$ws.Cells.Item(326, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(326, 3).Value = "La Araucanía"
$ws.Cells.Item(326, 4).Value = 45215
$ws.Cells.Item(326, 5).Value = 9
$ws.Cells.Item(326, 6).Value = 100112039
$ws.Cells.Item(326, 7).Value = "Ciboulette"
$ws.Cells.Item(326, 8).Value = "Sin especificar"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 125
$ws.Cells.Item(326, 11).Value = 6000
$ws.Cells.Item(326, 12).Value = 6000
$ws.Cells.Item(326, 13).Value = 6000
$ws.Cells.Item(326, 14).Value = "$/docena de atados"
$ws.Cells.Item(326, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(326, 16).Value = 2000
$ws.Cells.Item(326, 17).Value = 3
$ws.Cells.Item(326, 18).Value = "Hortaliza"
